$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the login email used in the "profile" row while keeping the
# password value the same.
$ws.Range("A2").Value = "tpnqatest@gmail.com"
$ws.Range("B2").Value = "Admin@123"
